$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 572; existing rows 572-694 shift down to 574-696.
$ws.Rows("572:573").Insert()

# New row 572 (Primera) - weekly reading for Femacal de La Calera / Betarraga
$ws.Cells.Item(572,1).Value  = 3
$ws.Cells.Item(572,2).Value  = "Femacal de La Calera"
$ws.Cells.Item(572,3).Value  = "Coquimbo"
$ws.Cells.Item(572,4).Value  = 44711
$ws.Cells.Item(572,5).Value  = 5
$ws.Cells.Item(572,6).Value  = 100114014
$ws.Cells.Item(572,7).Value  = "Betarraga"
$ws.Cells.Item(572,8).Value  = "Sin especificar"
$ws.Cells.Item(572,9).Value  = "Primera"
$ws.Cells.Item(572,10).Value = 3650
$ws.Cells.Item(572,11).Value = 550
$ws.Cells.Item(572,12).Value = 600
$ws.Cells.Item(572,13).Value = 575
$ws.Cells.Item(572,14).Value = "`$/paquete 4 unidades"
$ws.Cells.Item(572,15).Value = "Provincia de Quillota"
$ws.Cells.Item(572,16).Value = 144
$ws.Cells.Item(572,17).Value = 4
$ws.Cells.Item(572,18).Value = "Hortaliza"

# New row 573 (Segunda) - weekly reading for Femacal de La Calera / Betarraga
$ws.Cells.Item(573,1).Value  = 3
$ws.Cells.Item(573,2).Value  = "Femacal de La Calera"
$ws.Cells.Item(573,3).Value  = "Coquimbo"
$ws.Cells.Item(573,4).Value  = 44711
$ws.Cells.Item(573,5).Value  = 5
$ws.Cells.Item(573,6).Value  = 100114014
$ws.Cells.Item(573,7).Value  = "Betarraga"
$ws.Cells.Item(573,8).Value  = "Sin especificar"
$ws.Cells.Item(573,9).Value  = "Segunda"
$ws.Cells.Item(573,10).Value = 3700
$ws.Cells.Item(573,11).Value = 400
$ws.Cells.Item(573,12).Value = 450
$ws.Cells.Item(573,13).Value = 424
$ws.Cells.Item(573,14).Value = "`$/paquete 4 unidades"
$ws.Cells.Item(573,15).Value = "Provincia de Quillota"
$ws.Cells.Item(573,16).Value = 106
$ws.Cells.Item(573,17).Value = 4
$ws.Cells.Item(573,18).Value = "Hortaliza"
